$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying TPM data was recomputed, and the two "* -> ECs" rows
# (rows 6 and 7, Target cluster = ECs) were dropped from the LR-pair
# table. Delete them first so everything below shifts up cleanly.
$ws.Rows("6:7").Delete()

# Row 2: Sending cluster ECs -> Prok1/Prokr1 -> Target cluster now FAPs,
# with refreshed NATMI metrics from the new TPM values.
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("G2").Value2 = 0.412312
$ws.Range("H2").Value2 = 1.236936
$ws.Range("I2").Value2 = 0.7859709944540746
$ws.Range("J2").Value2 = 0.7859709944540746
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 3.661215333333333
$ws.Range("N2").Value2 = 10.983646
$ws.Range("O2").Value2 = 0.9813286907532969
$ws.Range("P2").Value2 = 0.9813286907532969
$ws.Range("Q2").Value2 = 1.509563016517333
$ws.Range("R2").Value2 = 13.586067148656
$ws.Range("S2").Value2 = 0.7712958869576838
$ws.Range("T2").Value2 = 0.7712958869576838
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("G3").Value2 = 0.412312
$ws.Range("H3").Value2 = 1.236936
$ws.Range("I3").Value2 = 0.7859709944540746
$ws.Range("J3").Value2 = 0.7859709944540746
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.06966033333333334
$ws.Range("N3").Value2 = 0.208981
$ws.Range("O3").Value2 = 0.01867130924670321
$ws.Range("P3").Value2 = 0.01867130924670321
$ws.Range("Q3").Value2 = 0.02872179135733334
$ws.Range("R3").Value2 = 0.258496122216
$ws.Range("S3").Value2 = 0.01467510749639088
$ws.Range("T3").Value2 = 0.01467510749639088
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("G4").Value2 = 0.1122773333333333
$ws.Range("H4").Value2 = 0.336832
$ws.Range("I4").Value2 = 0.2140290055459255
$ws.Range("J4").Value2 = 0.2140290055459255
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 3.661215333333333
$ws.Range("N4").Value2 = 10.983646
$ws.Range("O4").Value2 = 0.9813286907532969
$ws.Range("P4").Value2 = 0.9813286907532969
$ws.Range("Q4").Value2 = 0.4110714943857778
$ws.Range("R4").Value2 = 3.699643449472001
$ws.Range("S4").Value2 = 0.2100328037956131
$ws.Range("T4").Value2 = 0.2100328037956131
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("I5").Value2 = 0.2140290055459255
$ws.Range("J5").Value2 = 0.2140290055459255
$ws.Range("M5").Value2 = 0.06966033333333334
$ws.Range("N5").Value2 = 0.208981
$ws.Range("O5").Value2 = 0.01867130924670321
$ws.Range("P5").Value2 = 0.01867130924670321
$ws.Range("Q5").Value2 = 0.007821276465777779
$ws.Range("R5").Value2 = 0.070391488192
$ws.Range("S5").Value2 = 0.003996201750312331
$ws.Range("T5").Value2 = 0.003996201750312331
Write-Output "done"
